# Update bases das ligas - 17-05-2024 13:59
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 117 (cyclic rotation: gets former row 118 data) ---
$ws.Range("B117").Value = 7013702
$ws.Range("E117").Value = "Defensor Sporting"
$ws.Range("F117").Value = "Danubio"
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 2
$ws.Range("I117").Value = "A"
$ws.Range("J117").Value = 1.8
$ws.Range("K117").Value = 3.6
$ws.Range("L117").Value = 4.2
$ws.Range("M117").Value = 1.8
$ws.Range("N117").Value = 3.6
$ws.Range("O117").Value = 4.2
$ws.Range("P117").Value = -0.75
$ws.Range("Q117").Value = 2.05
$ws.Range("R117").Value = 1.8
$ws.Range("S117").Value = 2.25
$ws.Range("T117").Value = 1.85
$ws.Range("U117").Value = 2
$ws.Range("V117").Value = -1
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 3.2
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = 0.8
$ws.Range("AA117").Value = -0.5
$ws.Range("AB117").Value = 0.5

# --- Row 118 (cyclic rotation: gets former row 120 data) ---
$ws.Range("B118").Value = 7013409
$ws.Range("E118").Value = "Nacional De Football"
$ws.Range("F118").Value = "Torque"
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 1
$ws.Range("I118").Value = "D"
$ws.Range("J118").Value = 1.666
$ws.Range("K118").Value = 3.9
$ws.Range("L118").Value = 4.5
$ws.Range("M118").Value = 1.615
$ws.Range("N118").Value = 4
$ws.Range("O118").Value = 4.75
$ws.Range("P118").Value = -0.75
$ws.Range("Q118").Value = 1.8
$ws.Range("R118").Value = 2.05
$ws.Range("S118").Value = 2.75
$ws.Range("T118").Value = 1.95
$ws.Range("U118").Value = 1.9
$ws.Range("V118").Value = -1
$ws.Range("W118").Value = 3
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 1.05
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = 0.8999999999999999

# --- Row 120 (cyclic rotation: gets former row 117 data) ---
$ws.Range("B120").Value = 7013886
$ws.Range("E120").Value = "Racing Club de Montevideo"
$ws.Range("F120").Value = "Cerro"
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1
$ws.Range("I120").Value = "A"
$ws.Range("J120").Value = 2.25
$ws.Range("K120").Value = 3.1
$ws.Range("L120").Value = 3.25
$ws.Range("M120").Value = 2.25
$ws.Range("N120").Value = 2.875
$ws.Range("O120").Value = 3.5
$ws.Range("P120").Value = -0.25
$ws.Range("Q120").Value = 1.95
$ws.Range("R120").Value = 1.9
$ws.Range("S120").Value = 2
$ws.Range("T120").Value = 1.925
$ws.Range("U120").Value = 1.925
$ws.Range("V120").Value = -1
$ws.Range("W120").Value = -1
$ws.Range("X120").Value = 2.5
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = 0.8999999999999999
$ws.Range("AA120").Value = -1
$ws.Range("AB120").Value = 0.925

# --- Row 218 ---
$ws.Range("M218").Value = 3
$ws.Range("N218").Value = 3.1
$ws.Range("O218").Value = 2.4
$ws.Range("P218").Value = 0.25
$ws.Range("Q218").Value = 1.775
$ws.Range("R218").Value = 2.1
$ws.Range("T218").Value = 1.925
$ws.Range("U218").Value = 1.925

# --- Row 219 ---
$ws.Range("M219").Value = 2.1
$ws.Range("N219").Value = 3.2
$ws.Range("O219").Value = 3.7
$ws.Range("P219").Value = -0.25
$ws.Range("Q219").Value = 1.8
$ws.Range("R219").Value = 2.05

# --- Row 222 ---
$ws.Range("M222").Value = 1.95
$ws.Range("N222").Value = 3.4
$ws.Range("O222").Value = 4
$ws.Range("P222").Value = -0.5
$ws.Range("Q222").Value = 1.975
$ws.Range("R222").Value = 1.875
$ws.Range("T222").Value = 2.05
$ws.Range("U222").Value = 1.8
